$wb = $excel.ActiveWorkbook

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7667.3335
$ws.Range("J70").Value = 6500
$ws.Range("L70").Value = 19500
$ws.Range("N70").Value = -20040

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 7667.3335
$ws.Range("J73").Value = 6500
$ws.Range("L73").Value = 19500
$ws.Range("N73").Value = -21372

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1210.1428
$ws.Range("J129").Value = 1335.6364
$ws.Range("L129").Value = 4006.9092
$ws.Range("N129").Value = -14006.9092

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1913.5227
$ws.Range("I132").Value = 1318.5385
$ws.Range("J132").Value = 6554.4
$ws.Range("K132").Value = 3955.6155
$ws.Range("L132").Value = 19663.2
$ws.Range("M132").Value = -1425.6155
$ws.Range("N132").Value = -24723.2

# Sheet ALC, row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 46446.555
$ws.Range("J136").Value = 46446.555
$ws.Range("L136").Value = 46446.555
$ws.Range("N136").Value = -56646.555

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2120.25
$ws.Range("I137").Value = 1540.2727
$ws.Range("K137").Value = 4620.8181
$ws.Range("M137").Value = -2070.8181

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2659.2
$ws.Range("J138").Value = 3436.0212
$ws.Range("L138").Value = 10308.0636
$ws.Range("N138").Value = -20588.0636

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12110.115
$ws.Range("I32").Value = 9478.914000000001
$ws.Range("J32").Value = 19740.6
$ws.Range("K32").Value = 9478.914000000001
$ws.Range("L32").Value = 19740.6
$ws.Range("M32").Value = -9191.914000000001
$ws.Range("N32").Value = -20314.6

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 204083.34
$ws.Range("I61").Value = 4523.909
$ws.Range("K61").Value = 4523.909
$ws.Range("M61").Value = -4311.909

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1844.6666
$ws.Range("I74").Value = 1443.0714
$ws.Range("J74").Value = 2277.1538
$ws.Range("K74").Value = 1443.0714
$ws.Range("L74").Value = 2277.1538
$ws.Range("M74").Value = -569.0714
$ws.Range("N74").Value = -4025.1538

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1844.6666
$ws.Range("I77").Value = 1443.0714
$ws.Range("J77").Value = 2277.1538
$ws.Range("K77").Value = 7215.357
$ws.Range("L77").Value = 11385.769
$ws.Range("M77").Value = -2847.357
$ws.Range("N77").Value = -20121.769

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1451465.8
$ws.Range("I132").Value = 1373.5962
$ws.Range("J132").Value = 5887042
$ws.Range("K132").Value = 4120.7886
$ws.Range("L132").Value = 17661126
$ws.Range("M132").Value = -1590.7886
$ws.Range("N132").Value = -17666186

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 204083.34
$ws.Range("I136").Value = 4523.909
$ws.Range("K136").Value = 13571.727
$ws.Range("M136").Value = -11021.727

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 466.66666
$ws.Range("I7").Value = 600
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -487
$ws.Range("N7").Value = -426

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 50
$ws.Range("N22").ClearContents()

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 295327.5
$ws.Range("I31").Value = 1521.6129
$ws.Range("J31").Value = 709326.7
$ws.Range("K31").Value = 1521.6129
$ws.Range("L31").Value = 709326.7
$ws.Range("M31").Value = -1226.6129
$ws.Range("N31").Value = -709916.7

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 295327.5
$ws.Range("I34").Value = 1521.6129
$ws.Range("J34").Value = 709326.7
$ws.Range("K34").Value = 1521.6129
$ws.Range("L34").Value = 709326.7
$ws.Range("M34").Value = -1319.6129
$ws.Range("N34").Value = -709730.7

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 306601.5
$ws.Range("I134").Value = 3978.0908
$ws.Range("J134").Value = 911848.25
$ws.Range("K134").Value = 11934.2724
$ws.Range("L134").Value = 2735544.75
$ws.Range("M134").Value = -9399.2724
$ws.Range("N134").Value = -2740614.75

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6475.2383
$ws.Range("I5").Value = 9523.454
$ws.Range("J5").Value = 3122.2
$ws.Range("K5").Value = 28570.362
$ws.Range("L5").Value = 9366.599999999999
$ws.Range("M5").Value = -28458.362
$ws.Range("N5").Value = -9590.599999999999

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 6475.2383
$ws.Range("I135").Value = 9523.454
$ws.Range("J135").Value = 3122.2
$ws.Range("K135").Value = 85711.086
$ws.Range("L135").Value = 28099.8
$ws.Range("M135").Value = -83176.086
$ws.Range("N135").Value = -33169.8

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2430.6487
$ws.Range("I102").Value = 2229.348
$ws.Range("J102").Value = 2761.3572
$ws.Range("K102").Value = 2229.348
$ws.Range("L102").Value = 2761.3572
$ws.Range("M102").Value = -607.348
$ws.Range("N102").Value = -6005.3572

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3182.0657
$ws.Range("I132").Value = 2892.0698
$ws.Range("J132").Value = 3874.8333
$ws.Range("K132").Value = 8676.2094
$ws.Range("L132").Value = 11624.4999
$ws.Range("M132").Value = -6146.2094
$ws.Range("N132").Value = -16684.4999

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8338841.5
$ws.Range("I132").Value = 9809876
$ws.Range("J132").Value = 2983.1667
$ws.Range("K132").Value = 29429628
$ws.Range("L132").Value = 8949.500100000001
$ws.Range("M132").Value = -29427098
$ws.Range("N132").Value = -14009.5001

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 152514.67
$ws.Range("I62").Value = 152514.67
$ws.Range("K62").Value = 152514.67
$ws.Range("M62").Value = -151890.67

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 152514.67
$ws.Range("I65").Value = 152514.67
$ws.Range("K65").Value = 762573.3500000001
$ws.Range("M65").Value = -759453.3500000001

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 804
$ws.Range("I122").Value = 804
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2412
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 38
$ws.Range("N122").ClearContents()

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1869.881
$ws.Range("I132").Value = 1325.2413
$ws.Range("J132").Value = 3084.8462
$ws.Range("K132").Value = 3975.7239
$ws.Range("L132").Value = 9254.5386
$ws.Range("M132").Value = -1445.7239
$ws.Range("N132").Value = -14314.5386

# Sheet WVR, row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 44500
$ws.Range("J133").Value = 44500
$ws.Range("L133").Value = 44500
$ws.Range("N133").Value = -54620
